$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 29900
$ws.Cells.Item(21, 10).Value = 29900
$ws.Cells.Item(21, 12).Value = 29900
$ws.Cells.Item(21, 14).Value = -30836
$ws.Cells.Item(23, 8).Value = 29900
$ws.Cells.Item(23, 10).Value = 29900
$ws.Cells.Item(23, 12).Value = 29900
$ws.Cells.Item(23, 14).Value = -30368
$ws.Cells.Item(53, 8).Value = 220.2
$ws.Cells.Item(53, 9).Value = 76.5
$ws.Cells.Item(53, 10).Value = 795
$ws.Cells.Item(53, 11).Value = 76.5
$ws.Cells.Item(53, 12).Value = 795
$ws.Cells.Item(53, 13).Value = 560.5
$ws.Cells.Item(53, 14).Value = -2069
$ws.Cells.Item(58, 8).Value = 6980.6665
$ws.Cells.Item(58, 10).Value = 9998.5
$ws.Cells.Item(58, 12).Value = 29995.5
$ws.Cells.Item(58, 14).Value = -30295.5
$ws.Cells.Item(70, 8).Value = 1519.4
$ws.Cells.Item(70, 9).Value = 1534
$ws.Cells.Item(70, 11).Value = 4602
$ws.Cells.Item(70, 13).Value = -4332
$ws.Cells.Item(73, 8).Value = 1519.4
$ws.Cells.Item(73, 9).Value = 1534
$ws.Cells.Item(73, 11).Value = 4602
$ws.Cells.Item(73, 13).Value = -3666
$ws.Cells.Item(98, 8).Value = 819.2
$ws.Cells.Item(98, 9).Value = 819.2
$ws.Cells.Item(98, 11).Value = 819.2
$ws.Cells.Item(98, 13).Value = 678.8
$ws.Cells.Item(107, 8).Value = 1109
$ws.Cells.Item(107, 9).Value = 1078.9231
$ws.Cells.Item(107, 11).Value = 1078.9231
$ws.Cells.Item(107, 13).Value = 841.0769
$ws.Cells.Item(112, 8).Value = 1250.8462
$ws.Cells.Item(112, 10).Value = 1250.8462
$ws.Cells.Item(112, 12).Value = 3752.5386
$ws.Cells.Item(112, 14).Value = -5968.5386
$ws.Cells.Item(121, 8).Value = 1133.8334
$ws.Cells.Item(121, 10).Value = 1133.8334
$ws.Cells.Item(121, 12).Value = 3401.5002
$ws.Cells.Item(121, 14).Value = -6895.5002
$ws.Cells.Item(122, 8).Value = 819.2
$ws.Cells.Item(122, 9).Value = 819.2
$ws.Cells.Item(122, 11).Value = 2457.6
$ws.Cells.Item(122, 13).Value = -7.600000000000364
$ws.Cells.Item(125, 8).Value = 6960.3335
$ws.Cells.Item(125, 9).Value = 5945.5
$ws.Cells.Item(125, 10).Value = 8990
$ws.Cells.Item(125, 11).Value = 53509.5
$ws.Cells.Item(125, 12).Value = 80910
$ws.Cells.Item(125, 13).Value = -51049.5
$ws.Cells.Item(125, 14).Value = -85830
$ws.Cells.Item(131, 8).Value = 9800.200000000001
$ws.Cells.Item(131, 9).Value = 9974
$ws.Cells.Item(131, 11).Value = 29922
$ws.Cells.Item(131, 13).Value = -24882
$ws.Cells.Item(132, 8).Value = 1000.7143
$ws.Cells.Item(132, 9).Value = 1012.2
$ws.Cells.Item(132, 11).Value = 3036.6
$ws.Cells.Item(132, 13).Value = -506.6000000000004
$ws.Cells.Item(137, 8).Value = 3024.75
$ws.Cells.Item(137, 10).Value = 4250
$ws.Cells.Item(137, 12).Value = 12750
$ws.Cells.Item(137, 14).Value = -17850
$ws.Cells.Item(138, 8).Value = 8051.6875
$ws.Cells.Item(138, 10).Value = 8188.467
$ws.Cells.Item(138, 12).Value = 24565.401
$ws.Cells.Item(138, 14).Value = -34845.401
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2483.5
$ws.Cells.Item(45, 9).Value = 2163
$ws.Cells.Item(45, 10).Value = 3231.3333
$ws.Cells.Item(45, 11).Value = 2163
$ws.Cells.Item(45, 12).Value = 3231.3333
$ws.Cells.Item(45, 13).Value = -1786
$ws.Cells.Item(45, 14).Value = -3985.3333
$ws.Cells.Item(110, 8).Value = 4303.25
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(99, 8).Value = 1500
$ws.Cells.Item(99, 9).Value = 1500
$ws.Cells.Item(99, 11).Value = 1500
$ws.Cells.Item(99, 13).Value = -2
$ws.Cells.Item(107, 8).Value = 4608.6
$ws.Cells.Item(107, 9).Value = 4631.9443
$ws.Cells.Item(107, 10).Value = 4398.5
$ws.Cells.Item(107, 11).Value = 4631.9443
$ws.Cells.Item(107, 12).Value = 4398.5
$ws.Cells.Item(107, 13).Value = -2711.9443
$ws.Cells.Item(107, 14).Value = -8238.5
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 4595.7
$ws.Cells.Item(134, 9).Value = 4866.7144
$ws.Cells.Item(134, 10).Value = 3963.3333
$ws.Cells.Item(134, 11).Value = 14600.1432
$ws.Cells.Item(134, 12).Value = 11889.9999
$ws.Cells.Item(134, 13).Value = -12065.1432
$ws.Cells.Item(134, 14).Value = -16959.9999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6824.25
$ws.Cells.Item(31, 9).Value = 1783.8572
$ws.Cells.Item(31, 10).Value = 10744.556
$ws.Cells.Item(31, 11).Value = 1783.8572
$ws.Cells.Item(31, 12).Value = 10744.556
$ws.Cells.Item(31, 13).Value = -1488.8572
$ws.Cells.Item(31, 14).Value = -11334.556
$ws.Cells.Item(34, 8).Value = 6824.25
$ws.Cells.Item(34, 9).Value = 1783.8572
$ws.Cells.Item(34, 10).Value = 10744.556
$ws.Cells.Item(34, 11).Value = 1783.8572
$ws.Cells.Item(34, 12).Value = 10744.556
$ws.Cells.Item(34, 13).Value = -1581.8572
$ws.Cells.Item(34, 14).Value = -11148.556
$ws.Cells.Item(134, 8).Value = 2793.8462
$ws.Cells.Item(134, 9).Value = 2574.6365
$ws.Cells.Item(134, 11).Value = 7723.9095
$ws.Cells.Item(134, 13).Value = -5188.9095
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 20
$ws.Cells.Item(23, 9).Value = 20
$ws.Cells.Item(23, 11).Value = 60
$ws.Cells.Item(23, 13).Value = 175
$ws.Cells.Item(34, 8).Value = 1433.3334
$ws.Cells.Item(34, 9).Value = 1433.3334
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 4300.0002
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -4216.0002
$ws.Cells.Item(34, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 767.46155
$ws.Cells.Item(113, 10).Value = 891.5
$ws.Cells.Item(113, 12).Value = 2674.5
$ws.Cells.Item(113, 14).Value = -7014.5
$ws.Cells.Item(132, 8).Value = 4149.25
$ws.Cells.Item(132, 9).Value = 2056.5715
$ws.Cells.Item(132, 10).Value = 7079
$ws.Cells.Item(132, 11).Value = 18509.1435
$ws.Cells.Item(132, 12).Value = 63711
$ws.Cells.Item(132, 13).Value = -15979.1435
$ws.Cells.Item(132, 14).Value = -68771
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3144.4243
$ws.Cells.Item(132, 9).Value = 2630.4
$ws.Cells.Item(132, 11).Value = 7891.200000000001
$ws.Cells.Item(132, 13).Value = -5361.200000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(136, 8).Value = 34904.94
$ws.Cells.Item(136, 10).Value = 49089.547
$ws.Cells.Item(136, 12).Value = 147268.641
$ws.Cells.Item(136, 14).Value = -152368.641
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 3485.75
$ws.Cells.Item(96, 9).Value = 3485.75
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 3485.75
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -2112.75
$ws.Cells.Item(96, 14).ClearContents()
$ws.Cells.Item(133, 8).Value = 68000
$ws.Cells.Item(133, 10).Value = 68000
$ws.Cells.Item(133, 12).Value = 68000
$ws.Cells.Item(133, 14).Value = -78120
